$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.377.47"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "1.846.93"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.27"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6259"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9995"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07494"
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2899"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.39"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07736"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "1.846.20"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.995"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6801"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001044"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.12"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").Value = "2.105.95"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.178"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "29.407.39"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.45"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.33"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.454"
$ws.Range("E23").Value = "  -0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.59"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1376"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.405"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06454"
$ws.Range("E29").Value = "  +15.23%  "
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.477"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.063"
$ws.Range("E33").Value = "  -0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.827"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6993"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.578"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "1.262.84"
$ws.Range("E38").Value = "  +2.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.831"
$ws.Range("E39").Value = "  +4.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01829"
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.608"
$ws.Range("E41").Value = "  +3.46%  "
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9988"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "2.010.50"
$ws.Range("E44").Value = "  -18.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.39"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.29"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.757"
$ws.Range("E47").Value = "  +5.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.082"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1175"
$ws.Range("E49").Value = "  +3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.054"
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3948"
$ws.Range("E51").Value = "  -1.54%  "
